# Add Provantage price pull. Clean folder structure. Debug text.

$wb = $excel.ActiveWorkbook

$currencyFmt = "$#,##0.00_);[Red]($#,##0.00)"

# ---------------------------------------------------------------
# Newegg: C6 "See price in cart" -> numeric price; C7 price bump
# (done first so the now-unused "See price in cart" shared string
# slot is reclaimed by the next new strings written below)
# ---------------------------------------------------------------
$newegg = $wb.Worksheets.Item("Newegg")
$newegg.Range("C6").Value = 345.99
$newegg.Range("C6").NumberFormat = $currencyFmt
$newegg.Range("C7").Value = 1249

# ---------------------------------------------------------------
# New Provantage sheet + its pulled prices
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$provantage = $wb.Worksheets.Add($null, $lastSheet)
$provantage.Name = "Provantage"

$provantage.Range("B2").Value = "NECJ044"
$provantage.Range("C2").Value = 106.25

$provantage.Range("B3").Value = "AXM9718"
$provantage.Range("C3").Value = 57.21

$provantage.Range("B4").Value = "IOG90EM"
$provantage.Range("C4").Value = 21.17

$provantage.Range("B5").Value = "OCZT0WJ"
$provantage.Range("C5").Value = 51.2

# ---------------------------------------------------------------
# Zones: refreshed product link + new price pulled for it
# ---------------------------------------------------------------
$zones = $wb.Worksheets.Item("Zones")
$zones.Range("B2").Value = "http://www.zones.com/site/product/index.html?id=001598026"
$zones.Range("C2").Value = 34.99
$zones.Range("C2").NumberFormat = $currencyFmt

# ---------------------------------------------------------------
# Insight: refreshed part number + new price pulled for it
# ---------------------------------------------------------------
$insight = $wb.Worksheets.Item("Insight")
$insight.Range("B2").Value = "C7C95AW#ABA"
$insight.Range("C2").Value = 714.99

# ---------------------------------------------------------------
# Selections, restored/updated to match the refreshed state
# (order matters: last one selected becomes the active sheet/tab)
# ---------------------------------------------------------------
[void]$wb.Worksheets.Item("CDW").Range("C2:C4").Select()
[void]$wb.Worksheets.Item("PC Connections").Range("C2:C3").Select()
[void]$wb.Worksheets.Item("PCM").Range("C2:C3").Select()
[void]$insight.Range("B2").Select()
[void]$zones.Range("B2").Select()
[void]$newegg.Range("C2:C7").Select()
[void]$provantage.Range("C2").Select()

Write-Output "done"
